# Auto-generated Excel COM-interop script applying the Tonberry_Profits diff.
# Updates recompute per-item profit/cost figures across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 111.181816
$ws.Range("I33").Value = 77.8
$ws.Range("K33").Value = 77.8
$ws.Range("M33").Value = 151.2
$ws.Range("H55").Value = 150.7
$ws.Range("I55").Value = 156.11111
$ws.Range("K55").Value = 156.11111
$ws.Range("M55").Value = 57.88889
$ws.Range("H58").Value = 1388
$ws.Range("J58").Value = 1388
$ws.Range("L58").Value = 4164
$ws.Range("N58").Value = -4464
$ws.Range("H62").Value = 1974.75
$ws.Range("I62").Value = 1966.3334
$ws.Range("K62").Value = 1966.3334
$ws.Range("M62").Value = -1342.3334
$ws.Range("H65").Value = 1974.75
$ws.Range("I65").Value = 1966.3334
$ws.Range("K65").Value = 9831.666999999999
$ws.Range("M65").Value = -6711.666999999999
$ws.Range("H112").Value = 2175
$ws.Range("J112").Value = 2221.739
$ws.Range("L112").Value = 6665.217000000001
$ws.Range("N112").Value = -8881.217000000001
$ws.Range("H127").Value = 1846.6086
$ws.Range("I127").Value = 2143.0908
$ws.Range("J127").Value = 1574.8334
$ws.Range("K127").Value = 6429.2724
$ws.Range("L127").Value = 4724.5002
$ws.Range("M127").Value = -1469.2724
$ws.Range("N127").Value = -14644.5002
$ws.Range("H131").Value = 2291.3447
$ws.Range("I131").Value = 749.46155
$ws.Range("K131").Value = 2248.38465
$ws.Range("M131").Value = 2791.61535
$ws.Range("H132").Value = 7093577.5
$ws.Range("I132").Value = 8334634
$ws.Range("J132").Value = 1826
$ws.Range("K132").Value = 25003902
$ws.Range("L132").Value = 5478
$ws.Range("M132").Value = -25001372
$ws.Range("N132").Value = -10538
$ws.Range("H135").Value = 649.0454999999999
$ws.Range("I135").Value = 618.95
$ws.Range("K135").Value = 5570.55
$ws.Range("M135").Value = -3035.55
$ws.Range("H137").Value = 1176.7878
$ws.Range("I137").Value = 1009.44446
$ws.Range("J137").Value = 1377.6
$ws.Range("K137").Value = 3028.33338
$ws.Range("L137").Value = 4132.799999999999
$ws.Range("M137").Value = -478.33338
$ws.Range("N137").Value = -9232.799999999999
$ws.Range("H138").Value = 1763.2347
$ws.Range("I138").Value = 1450.6492
$ws.Range("J138").Value = 2197.805
$ws.Range("K138").Value = 4351.9476
$ws.Range("L138").Value = 6593.414999999999
$ws.Range("M138").Value = 788.0523999999996
$ws.Range("N138").Value = -16873.415
$ws.Range("H141").Value = 1122224.8
$ws.Range("I141").Value = 1334207.8
$ws.Range("J141").Value = 9314
$ws.Range("K141").Value = 4002623.4
$ws.Range("L141").Value = 27942
$ws.Range("M141").Value = -3997443.4
$ws.Range("N141").Value = -38302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 327312.4
$ws.Range("I2").Value = 427837.47
$ws.Range("K2").Value = 427837.47
$ws.Range("M2").Value = -427724.47
$ws.Range("H32").Value = 3527.2666
$ws.Range("I32").Value = 3285.415
$ws.Range("K32").Value = 3285.415
$ws.Range("M32").Value = -2998.415
$ws.Range("H45").Value = 1836.8948
$ws.Range("I45").Value = 1834.8889
$ws.Range("J45").Value = 1838.7
$ws.Range("K45").Value = 1834.8889
$ws.Range("L45").Value = 1838.7
$ws.Range("M45").Value = -1457.8889
$ws.Range("N45").Value = -2592.7
$ws.Range("H61").Value = 29413696
$ws.Range("I61").Value = 22728804
$ws.Range("J61").Value = 41669330
$ws.Range("K61").Value = 22728804
$ws.Range("L61").Value = 41669330
$ws.Range("M61").Value = -22728592
$ws.Range("N61").Value = -41669754
$ws.Range("H74").Value = 1310.9117
$ws.Range("I74").Value = 1098.4348
$ws.Range("J74").Value = 1755.1818
$ws.Range("K74").Value = 1098.4348
$ws.Range("L74").Value = 1755.1818
$ws.Range("M74").Value = -224.4348
$ws.Range("N74").Value = -3503.1818
$ws.Range("H77").Value = 1310.9117
$ws.Range("I77").Value = 1098.4348
$ws.Range("J77").Value = 1755.1818
$ws.Range("K77").Value = 5492.174
$ws.Range("L77").Value = 8775.909
$ws.Range("M77").Value = -1124.174
$ws.Range("N77").Value = -17511.909
$ws.Range("H110").Value = 3070.2778
$ws.Range("I110").Value = 2202.5
$ws.Range("K110").Value = 2202.5
$ws.Range("M110").Value = -157.5
$ws.Range("H116").Value = 327312.4
$ws.Range("I116").Value = 427837.47
$ws.Range("K116").Value = 427837.47
$ws.Range("M116").Value = -425543.47
$ws.Range("H122").Value = 1956.6296
$ws.Range("I122").Value = 1520.579
$ws.Range("K122").Value = 4561.737
$ws.Range("M122").Value = -2111.737
$ws.Range("H132").Value = 1409.7347
$ws.Range("I132").Value = 1127.2894
$ws.Range("J132").Value = 2385.4546
$ws.Range("K132").Value = 3381.8682
$ws.Range("L132").Value = 7156.3638
$ws.Range("M132").Value = -851.8681999999999
$ws.Range("N132").Value = -12216.3638
$ws.Range("H136").Value = 29413696
$ws.Range("I136").Value = 22728804
$ws.Range("J136").Value = 41669330
$ws.Range("K136").Value = 68186412
$ws.Range("L136").Value = 125007990
$ws.Range("M136").Value = -68183862
$ws.Range("N136").Value = -125013090
$ws.Range("H141").Value = 58429
$ws.Range("J141").Value = 58429
$ws.Range("L141").Value = 58429
$ws.Range("N141").Value = -68789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 327312.4
$ws.Range("I3").Value = 427837.47
$ws.Range("K3").Value = 427837.47
$ws.Range("M3").Value = -427723.47
$ws.Range("H20").Value = 2421
$ws.Range("I20").Value = 1984
$ws.Range("K20").Value = 1984
$ws.Range("M20").Value = -1737
$ws.Range("H94").Value = 502
$ws.Range("I94").Value = 569.6667
$ws.Range("K94").Value = 569.6667
$ws.Range("M94").Value = -118.6667
$ws.Range("H105").Value = 2225.9678
$ws.Range("I105").Value = 2246.3845
$ws.Range("J105").Value = 2119.8
$ws.Range("K105").Value = 2246.3845
$ws.Range("L105").Value = 2119.8
$ws.Range("M105").Value = -499.3845000000001
$ws.Range("N105").Value = -5613.8
$ws.Range("H107").Value = 1899
$ws.Range("I107").Value = 1899
$ws.Range("K107").Value = 1899
$ws.Range("M107").Value = 21
$ws.Range("H134").Value = 4327.4243
$ws.Range("I134").Value = 4569.871
$ws.Range("K134").Value = 13709.613
$ws.Range("M134").Value = -11174.613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1137.4375
$ws.Range("I22").Value = 424.875
$ws.Range("J22").Value = 1850
$ws.Range("K22").Value = 424.875
$ws.Range("L22").Value = 1850
$ws.Range("M22").Value = -74.875
$ws.Range("N22").Value = -2550
$ws.Range("H31").Value = 3573422.5
$ws.Range("I31").Value = 5495972.5
$ws.Range("K31").Value = 5495972.5
$ws.Range("M31").Value = -5495677.5
$ws.Range("H34").Value = 3573422.5
$ws.Range("I34").Value = 5495972.5
$ws.Range("K34").Value = 5495972.5
$ws.Range("M34").Value = -5495770.5
$ws.Range("H58").Value = 1891385.2
$ws.Range("I58").Value = 2288758.5
$ws.Range("K58").Value = 2288758.5
$ws.Range("M58").Value = -2288555.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H99").Value = 1695.8
$ws.Range("I99").Value = 1619.75
$ws.Range("K99").Value = 1619.75
$ws.Range("M99").Value = -121.75
$ws.Range("H107").Value = 2236.6667
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 3446.5454
$ws.Range("I122").Value = 3122.6155
$ws.Range("J122").Value = 3914.4443
$ws.Range("K122").Value = 9367.8465
$ws.Range("L122").Value = 11743.3329
$ws.Range("M122").Value = -6917.8465
$ws.Range("N122").Value = -16643.3329
$ws.Range("H126").Value = 1695.8
$ws.Range("I126").Value = 1619.75
$ws.Range("K126").Value = 4859.25
$ws.Range("M126").Value = -2389.25
$ws.Range("H132").Value = 1701.1765
$ws.Range("I132").Value = 1175.9584
$ws.Range("J132").Value = 2961.7
$ws.Range("K132").Value = 3527.8752
$ws.Range("L132").Value = 8885.099999999999
$ws.Range("M132").Value = -997.8751999999999
$ws.Range("N132").Value = -13945.1
$ws.Range("H134").Value = 1674.9608
$ws.Range("I134").Value = 1550.925
$ws.Range("J134").Value = 2126
$ws.Range("K134").Value = 4652.775
$ws.Range("L134").Value = 6378
$ws.Range("M134").Value = -2117.775
$ws.Range("N134").Value = -11448
$ws.Range("H136").Value = 1891385.2
$ws.Range("I136").Value = 2288758.5
$ws.Range("K136").Value = 6866275.5
$ws.Range("M136").Value = -6863725.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 932.5
$ws.Range("I3").Value = 932.5
$ws.Range("K3").Value = 2797.5
$ws.Range("M3").Value = -2685.5
$ws.Range("H26").Value = 1102.75
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 3000
$ws.Range("N26").Value = -3576
$ws.Range("H43").Value = 3250
$ws.Range("J43").Value = 3250
$ws.Range("L43").Value = 9750
$ws.Range("N43").Value = -9978
$ws.Range("H56").Value = 7206.4375
$ws.Range("I56").Value = 7206.4375
$ws.Range("K56").Value = 7206.4375
$ws.Range("M56").Value = -6676.4375
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 873.0857
$ws.Range("I122").Value = 512.7143
$ws.Range("K122").Value = 4614.428699999999
$ws.Range("M122").Value = -2164.428699999999
$ws.Range("H131").Value = 871.08
$ws.Range("I131").Value = 505.42856
$ws.Range("J131").Value = 898.6022
$ws.Range("K131").Value = 1516.28568
$ws.Range("L131").Value = 2695.8066
$ws.Range("M131").Value = 3523.71432
$ws.Range("N131").Value = -12775.8066
$ws.Range("H139").Value = 5761.625
$ws.Range("I139").Value = 6103.591
$ws.Range("K139").Value = 18310.773
$ws.Range("M139").Value = -13170.773

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2977.182
$ws.Range("I102").Value = 3024.9
$ws.Range("K102").Value = 3024.9
$ws.Range("M102").Value = -1402.9
$ws.Range("H122").Value = 1627
$ws.Range("I122").Value = 1628.3334
$ws.Range("K122").Value = 4885.0002
$ws.Range("M122").Value = -2435.0002
$ws.Range("H126").Value = 2573237
$ws.Range("J126").Value = 61155.59
$ws.Range("L126").Value = 183466.77
$ws.Range("N126").Value = -188406.77
$ws.Range("H132").Value = 856364.4399999999
$ws.Range("I132").Value = 962830.0600000001
$ws.Range("J132").Value = 4639.2
$ws.Range("K132").Value = 2888490.18
$ws.Range("L132").Value = 13917.6
$ws.Range("M132").Value = -2885960.18
$ws.Range("N132").Value = -18977.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2876.5293
$ws.Range("I7").Value = 2262.75
$ws.Range("K7").Value = 2262.75
$ws.Range("M7").Value = -2150.75
$ws.Range("H22").Value = 3082.5715
$ws.Range("I22").Value = 5424.5
$ws.Range("J22").Value = 2145.8
$ws.Range("K22").Value = 5424.5
$ws.Range("L22").Value = 2145.8
$ws.Range("M22").Value = -5129.5
$ws.Range("N22").Value = -2735.8
$ws.Range("H27").Value = 3082.5715
$ws.Range("I27").Value = 5424.5
$ws.Range("J27").Value = 2145.8
$ws.Range("K27").Value = 5424.5
$ws.Range("L27").Value = 2145.8
$ws.Range("M27").Value = -5317.5
$ws.Range("N27").Value = -2359.8
$ws.Range("H40").Value = 9186.611000000001
$ws.Range("I40").Value = 9881.923000000001
$ws.Range("K40").Value = 9881.923000000001
$ws.Range("M40").Value = -9745.923000000001
$ws.Range("H122").Value = 4370.7144
$ws.Range("J122").Value = 5399
$ws.Range("L122").Value = 16197
$ws.Range("N122").Value = -21097
$ws.Range("H126").Value = 2876.5293
$ws.Range("I126").Value = 2262.75
$ws.Range("K126").Value = 6788.25
$ws.Range("M126").Value = -4318.25
$ws.Range("H132").Value = 1921.7587
$ws.Range("I132").Value = 1334.174
$ws.Range("K132").Value = 4002.522
$ws.Range("M132").Value = -1472.522
$ws.Range("H136").Value = 2153.0425
$ws.Range("I136").Value = 1293.5946
$ws.Range("K136").Value = 3880.7838
$ws.Range("M136").Value = -1330.7838

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8990.532999999999
$ws.Range("I14").Value = 8993.333000000001
$ws.Range("J14").Value = 8989.833000000001
$ws.Range("K14").Value = 8993.333000000001
$ws.Range("L14").Value = 8989.833000000001
$ws.Range("M14").Value = -8825.333000000001
$ws.Range("N14").Value = -9325.833000000001
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680
$ws.Range("H113").Value = 667.6
$ws.Range("I113").Value = 334.5
$ws.Range("K113").Value = 1003.5
$ws.Range("M113").Value = 1166.5
$ws.Range("H122").Value = 88114.89
$ws.Range("I122").Value = 98879.25
$ws.Range("K122").Value = 296637.75
$ws.Range("M122").Value = -294187.75
$ws.Range("H132").Value = 1436.375
$ws.Range("I132").Value = 993.3461
$ws.Range("K132").Value = 2980.0383
$ws.Range("M132").Value = -450.0383000000002
$ws.Range("H136").Value = 1475.9796
$ws.Range("I136").Value = 1343.725
$ws.Range("K136").Value = 4031.175
$ws.Range("M136").Value = -1481.175

